$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "257.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.49%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.59%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.739"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-11.03%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05972"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.21%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.674"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.45%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8710"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.73%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9571"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.51%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1410"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.22%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07183"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.05%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03142"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.15%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09233"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.02%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001544"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.25%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006083"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.59%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.11%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.486"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.207"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.51%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.218"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.33%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.72%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03606"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.81%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1307"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.19%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.528"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.01%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04233"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.93%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1381"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.16%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.65%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004505"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-11.84%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.07%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001494"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-22.89%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03825"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.86%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006233"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "10.15%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1102"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.26%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.62%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01102"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.16%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005498"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.84%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08554"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-4.76%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002269"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "5.30%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
